$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header labels (drop the internal "_x" suffix / stray spaces
# left over from the pandas column-flattening step).
$ws.Range("B1").Value = "Canada_priceprice"
$ws.Range("D1").Value = "Canada_pointspoints"

# Column C only ever held the stray "_1" header (no data below it) -
# delete it so the "points" data in column D slides left into column C.
$ws.Columns.Item(3).Delete()

# Column C (formerly D) gets widened to fit the merged header text;
# columns A/B keep their original widths.
$ws.Columns.Item(3).ColumnWidth = 20.83

Write-Host "Done"
